# Update latest output (run 51)

$wb = $excel.ActiveWorkbook

# --- Sheet "Schedule" ---
$schedule = $wb.Worksheets.Item("Schedule")
$schedule.Range("E2").Value = 853.8373349999999
$schedule.Range("F2").Value = 14.11768080357143

# --- Sheet "Detailed" ---
$detailed = $wb.Worksheets.Item("Detailed")

$detailed.Range("B7").Value = 36.06
$detailed.Range("B8").Value = 36.06

$detailed.Range("B9").Value = 57.09
$detailed.Range("C9").Value = "historical"

$detailed.Range("B10").Value = 57.09
$detailed.Range("C10").Value = "historical"

$detailed.Range("B11").Value = 68.5436

$detailed.Range("B18").Value = 36.06011
$detailed.Range("B19").Value = 36.06011
$detailed.Range("B21").Value = 36.06011
$detailed.Range("B22").Value = 46.3948

$detailed.Range("B26").Value = 36.06

$detailed.Range("B29").Value = 36.06029
$detailed.Range("B30").Value = 36.06029

$detailed.Range("B35").Value = 27.14153
$detailed.Range("B36").Value = -0.35021
$detailed.Range("B37").Value = -3.01616
$detailed.Range("B38").Value = -2.85253
$detailed.Range("B39").Value = -2.79946

$detailed.Range("B41").Value = 9.456770000000001
$detailed.Range("B42").Value = 9.784940000000001
$detailed.Range("B43").Value = 26.13075
$detailed.Range("B44").Value = 9.433719999999999
$detailed.Range("B45").Value = 9.75508
